$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"
